$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AO1").Value = "Answer 6"

$values = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 41).Value = [bool]($values[$i])
}

$ws.Range("AE1").Select()
$ws.Range("AP10").Select()
